# AGA206 Assessment 2 Checklist - mark "Rotating Wall or Door" (J24) and
# "Moving Hazard" (J32) as completed by checking their linked checkboxes.
# Checking the box sets the linked cell (J24/J32) to TRUE; the Status
# ("To Be Done" -> "Done"), per-row point total, and the summary totals
# (D8, D9, K40) all recalculate automatically from existing formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

$ws.Range("J24").Value = $true
$ws.Range("J32").Value = $true

# Reflect the navigation state left behind after ticking the boxes.
[void]$ws.Range("C24").Select()
